$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:50:38.132919",
    "2021-10-05 10:50:38.132930",
    "2021-10-05 10:50:38.132934",
    "2021-10-05 10:50:38.132936",
    "2021-10-05 10:50:38.132939",
    "2021-10-05 10:50:38.132942",
    "2021-10-05 10:50:38.132945",
    "2021-10-05 10:50:38.132948",
    "2021-10-05 10:50:38.132951",
    "2021-10-05 10:50:38.132954",
    "2021-10-05 10:50:38.132956",
    "2021-10-05 10:50:38.132959",
    "2021-10-05 10:50:38.132962",
    "2021-10-05 10:50:38.132965",
    "2021-10-05 10:50:38.132967",
    "2021-10-05 10:50:38.132970",
    "2021-10-05 10:50:38.132973",
    "2021-10-05 10:50:38.132976",
    "2021-10-05 10:50:38.132979",
    "2021-10-05 10:50:38.132982",
    "2021-10-05 10:50:38.132984",
    "2021-10-05 10:50:38.132987",
    "2021-10-05 10:50:38.132989",
    "2021-10-05 10:50:38.132992",
    "2021-10-05 10:50:38.132995",
    "2021-10-05 10:50:38.132998",
    "2021-10-05 10:50:38.133000",
    "2021-10-05 10:50:38.133003",
    "2021-10-05 10:50:38.133006",
    "2021-10-05 10:50:38.133009",
    "2021-10-05 10:50:38.133011",
    "2021-10-05 10:50:38.133014",
    "2021-10-05 10:50:38.133017"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

Write-Output "Done"
